{"js": "// Replace each \"before\" division expression with its corresponding\n// \"after\" value, matching the unified diff exactly (25 cell edits,\n// each old value unique, each new value unique).\nconst replacements = [\n  [\"647\u00f75=\", \"796\u00f72=\"],\n  [\"779\u00f72=\", \"150\u00f75=\"],\n  [\"278\u00f73=\", \"493\u00f75=\"],\n  [\"660\u00f75=\", \"699\u00f74=\"],\n  [\"814\u00f76=\", \"737\u00f75=\"],\n  [\"381\u00f76=\", \"245\u00f77=\"],\n  [\"141\u00f78=\", \"531\u00f77=\"],\n  [\"389\u00f75=\", \"290\u00f75=\"],\n  [\"325\u00f73=\", \"172\u00f73=\"],\n  [\"544\u00f79=\", \"983\u00f78=\"],\n  [\"959\u00f79=\", \"418\u00f78=\"],\n  [\"479\u00f78=\", \"695\u00f75=\"],\n  [\"103\u00f78=\", \"672\u00f76=\"],\n  [\"796\u00f73=\", \"573\u00f76=\"],\n  [\"238\u00f77=\", \"506\u00f78=\"],\n  [\"274\u00f78=\", \"451\u00f79=\"],\n  [\"238\u00f73=\", \"569\u00f75=\"],\n  [\"984\u00f74=\", \"297\u00f78=\"],\n  [\"821\u00f74=\", \"402\u00f72=\"],\n  [\"246\u00f75=\", \"728\u00f75=\"],\n  [\"236\u00f74=\", \"301\u00f72=\"],\n  [\"139\u00f73=\", \"307\u00f79=\"],\n  [\"958\u00f77=\", \"265\u00f78=\"],\n  [\"866\u00f73=\", \"426\u00f77=\"],\n  [\"541\u00f76=\", \"104\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"before\" division expression with its corresponding\n# \"after\" value, matching the unified diff exactly (25 cell edits,\n# each old value unique, each new value unique).\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = [ordered]@{\n    \"647\u00f75=\" = \"796\u00f72=\"\n    \"779\u00f72=\" = \"150\u00f75=\"\n    \"278\u00f73=\" = \"493\u00f75=\"\n    \"660\u00f75=\" = \"699\u00f74=\"\n    \"814\u00f76=\" = \"737\u00f75=\"\n    \"381\u00f76=\" = \"245\u00f77=\"\n    \"141\u00f78=\" = \"531\u00f77=\"\n    \"389\u00f75=\" = \"290\u00f75=\"\n    \"325\u00f73=\" = \"172\u00f73=\"\n    \"544\u00f79=\" = \"983\u00f78=\"\n    \"959\u00f79=\" = \"418\u00f78=\"\n    \"479\u00f78=\" = \"695\u00f75=\"\n    \"103\u00f78=\" = \"672\u00f76=\"\n    \"796\u00f73=\" = \"573\u00f76=\"\n    \"238\u00f77=\" = \"506\u00f78=\"\n    \"274\u00f78=\" = \"451\u00f79=\"\n    \"238\u00f73=\" = \"569\u00f75=\"\n    \"984\u00f74=\" = \"297\u00f78=\"\n    \"821\u00f74=\" = \"402\u00f72=\"\n    \"246\u00f75=\" = \"728\u00f75=\"\n    \"236\u00f74=\" = \"301\u00f72=\"\n    \"139\u00f73=\" = \"307\u00f79=\"\n    \"958\u00f77=\" = \"265\u00f78=\"\n    \"866\u00f73=\" = \"426\u00f77=\"\n    \"541\u00f76=\" = \"104\u00f76=\"\n}\n\nforeach ($old in $replacements.Keys) {\n    $new = $replacements[$old]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $new, $wdReplaceAll)\n}\n"}
